# Capitalize the PHENOTYPE column (column A, rows 2-27) so that each
# phenotype label starts with an uppercase letter.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $text = [string]$cell.Value2
    if ($text.Length -gt 0) {
        $cell.Value = $text.Substring(0,1).ToUpper() + $text.Substring(1)
    }
}

# Update the view state to match: scrolled so row 11 is at the top,
# with cell A27 selected.
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("A27").Select()
